# The English resource file dropped four "Derivative" radio-button strings
# (strRadBackwardOne, strRadCentralFive, strRadCentralThree, strRadForwardOne)
# so the Greek translation table must drop the matching rows too.
#
# The worksheet is an alphabetically-sorted Key/Comment/English/Greek table
# (Excel Table "Tabla13") bound to an autoFilter + sort on column B. Deleting
# the obsolete rows lets Excel naturally shift everything below up, shrink
# the table/used-range from B2:E169 to B2:E165, and drop the now-unused
# shared strings on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the 4 rows to remove by their Key (column B) so the script is
# resilient even if row numbers ever shift before this code runs.
$keysToRemove = @("strRadBackwardOne", "strRadCentralFive", "strRadCentralThree", "strRadForwardOne")

$lastRow = $ws.Range("B2").End(4).Row  # xlDown = 4 -> last contiguous filled row in column B
$rowsToDelete = @()

for ($r = 2; $r -le $lastRow; $r++) {
    $key = $ws.Cells.Item($r, 2).Value2
    if ($keysToRemove -contains $key) {
        $rowsToDelete += $r
    }
}

# Delete from the bottom up so earlier row numbers stay valid as we go.
$sortedRows = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sortedRows) {
    $ws.Rows.Item($r).Delete()
}
